# Apply the "I0 and IF added" change:
#  - Add header "I0" in I1 and "IF" in J1 (bold/bordered header style, like H1)
#  - Fill column I (I0) and column J (IF) for rows 2..22
#  - Dimension grows from A1:H22 to A1:J22 (handled automatically by Excel)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the new headers.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Values for the new I0 / IF columns, keyed by row number.
$values = @{
    2  = @(1, 4)
    3  = @(1, 4)
    4  = @(1, 3)
    5  = @(1, 2)
    6  = @(1, 2)
    7  = @(1, 2)
    8  = @(1, 3)
    9  = @(1, 3)
    10 = @(1, 4)
    11 = @(1, 3)
    12 = @(1, 3)
    13 = @(6, 7)
    14 = @(1, 5)
    15 = @(1, 3)
    16 = @(1, 5)
    17 = @(1, 4)
    18 = @(5, 6)
    19 = @(7, 8)
    20 = @(6, 6)
    21 = @(1, 2)
    22 = @(3, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # column J
}
